$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $excel.GetType().FullName
Write-Host $wb.GetType().FullName
Write-Host $ws.GetType().FullName
$cell = $ws.Range("A1")
Write-Host $cell.GetType().FullName
$members = $cell | Get-Member
Write-Host $members
